$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the first sheet
$ws.Name = "Export as TSV"

# 2. Freeze the header row (split below row 1, freeze panes)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Add errorTitle / error messages to the existing data validations
$val = $ws.Range("I2:I1048576").Validation
$val.ErrorTitle = "Value must come from list"
$val.ErrorMessage = "Value must be one of: sequence."

$val = $ws.Range("J2:J1048576").Validation
$val.ErrorTitle = "Value must come from list"
$val.ErrorMessage = "Value must be one of: Slide-seq."

$val = $ws.Range("K2:K1048576").Validation
$val.ErrorTitle = "Value must come from list"
$val.ErrorMessage = "Value must be one of: RNA."

$val = $ws.Range("L2:L1048576").Validation
$val.ErrorTitle = "Not a boolean"
$val.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$val = $ws.Range("T2:T1048576").Validation
$val.ErrorTitle = "Not a boolean"
$val.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$val = $ws.Range("Z2:Z1048576").Validation
$val.ErrorTitle = "Not a number"
$val.ErrorMessage = "The values in this column must be numbers."

$val = $ws.Range("AA2:AA1048576").Validation
$val.ErrorTitle = "Value must come from list"
$val.ErrorMessage = "Value must be one of: ng."

$val = $ws.Range("AE2:AE1048576").Validation
$val.ErrorTitle = "Not a number"
$val.ErrorMessage = "The values in this column must be numbers."

$val = $ws.Range("AF2:AF1048576").Validation
$val.ErrorTitle = "Not a number"
$val.ErrorMessage = "The values in this column must be numbers."
